# Feria Lagunitas de Puerto Montt - Coliflor: add a new weekly record.
# A new row of data is inserted at row 67 (pushing the existing rows 67-189
# down to 68-190), and the new row is populated with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67; everything below (rows 67-189) shifts
# down by one (to rows 68-190), automatically carrying the former last row
# (189) into the new row 190.
$ws.Rows("67:67").Insert()

# Populate the newly-inserted (blank) row 67 with its data.
$ws.Range("A67").Value = 4
$ws.Range("B67").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C67").Value = "Los Lagos"
$ws.Range("D67").Value = 44469
$ws.Range("E67").Value = 10
$ws.Range("F67").Value = 100112008
$ws.Range("G67").Value = "Coliflor"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 500
$ws.Range("K67").Value = 1300
$ws.Range("L67").Value = 1300
$ws.Range("M67").Value = 1300
$ws.Range("N67").Value = "$/unidad"
$ws.Range("O67").Value = "Región Metropolitana"
$ws.Range("P67").Value = 1300
$ws.Range("Q67").Value = 1
$ws.Range("R67").Value = "Hortaliza"
